$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) retains text formatting so numeric-looking
# strings such as "24.70" or "0.00001060" are not coerced into numbers
# and do not lose formatting (matches original inlineStr cell content).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.368.61'
$ws.Range("E2").Value = '  -0.07%  '
$ws.Range("D3").Value = '1.846.96'
$ws.Range("E3").Value = '  -0.08%  '
$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '240.18'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").Value = '0.6277'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("D7").Value = '0.9993'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '0.07637'
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("E9").Value = '  -1.00%  '
$ws.Range("D10").Value = '24.70'
$ws.Range("E10").Value = '  +0.74%  '
$ws.Range("D11").Value = '0.07736'
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '5.030'
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '0.6779'
$ws.Range("E13").Value = '  -0.14%  '
$ws.Range("D14").Value = '0.00001060'
$ws.Range("E14").Value = '  -2.26%  '
$ws.Range("D16").Value = '6.153'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '29.410.41'
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("D18").Value = '227.13'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D20").Value = '0.9993'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '7.490'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = '158.21'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("D24").Value = '0.1383'
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").Value = '8.402'
$ws.Range("E25").Value = '  +0.33%  '
$ws.Range("D26").Value = '17.70'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").Value = '1.382'
$ws.Range("E27").Value = '  +5.04%  '
$ws.Range("D28").Value = '1.459'
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("D29").Value = '0.05606'
$ws.Range("D30").Value = '4.117'
$ws.Range("E30").Value = '  +0.17%  '
$ws.Range("D31").Value = '4.099'
$ws.Range("E31").Value = '  +1.52%  '
$ws.Range("D32").Value = '1.835'
$ws.Range("E32").Value = '  -0.69%  '
$ws.Range("D33").Value = '1.161'
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("D34").Value = '0.6939'
$ws.Range("E34").Value = '  -2.12%  '
$ws.Range("D35").Value = '2.576'
$ws.Range("E35").Value = '  -0.32%  '
$ws.Range("D36").Value = '0.01802'
$ws.Range("E36").Value = '  +0.26%  '
$ws.Range("D37").Value = '1.227.43'
$ws.Range("E37").Value = '  -0.42%  '
$ws.Range("D38").Value = '2.712'
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").Value = '6.385'
$ws.Range("E39").Value = '  -1.30%  '
$ws.Range("D40").Value = '0.9054'
$ws.Range("E40").Value = '  -0.30%  '
$ws.Range("D41").Value = '0.9993'
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = '101.47'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").Value = '66.02'
$ws.Range("D44").Value = '7.172'
$ws.Range("E44").Value = '  +0.25%  '
$ws.Range("D45").Value = '0.4007'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("D46").Value = '9.012'
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = '1.674'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("D48").Value = '0.1140'
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("D49").Value = '0.05702'
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("D50").Value = '0.4625'
$ws.Range("E51").Value = '  +0.26%  '
